$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SmartRules")
Write-Host $ws.Name
